# Add "Legal Vehicle" column to the Operation Masterlist report
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert a new column at C (shifts Date.. and everything after one column right,
#    and automatically grows the A1:C1 title merge to A1:D1 and the NCD merges)
$ws.Columns("C:C").Insert()

# 2. New header text for the inserted column
$ws.Range("C3").Value = "Legal Vehicle"

# 3. Column widths (ColumnWidth is in character units; stored xlsx width = ColumnWidth + 0.8333333333333334)
$ws.Columns("A:A").ColumnWidth = 10.166666666666666   # 11
$ws.Columns("B:C").ColumnWidth = 23.592447916666668   # 24.42578125
$ws.Columns("D:D").ColumnWidth = 9.877604166666666    # 10.7109375
$ws.Columns("E:E").ColumnWidth = 7.166666666666667    # 8
$ws.Columns("F:F").ColumnWidth = 8.736979166666666    # 9.5703125
$ws.Columns("G:G").ColumnWidth = 3.5924479166666665   # 4.42578125
$ws.Columns("H:H").ColumnWidth = 6.877604166666667    # 7.7109375
$ws.Columns("I:I").ColumnWidth = 20.022135416666668   # 20.85546875
$ws.Columns("J:J").ColumnWidth = 9.877604166666666    # 10.7109375
$ws.Columns("K:K").ColumnWidth = 7.877604166666667    # 8.7109375
$ws.Columns("L:L").ColumnWidth = 20.022135416666668   # 20.85546875
$ws.Columns("M:M").ColumnWidth = 9.877604166666666    # 10.7109375
$ws.Columns("N:N").ColumnWidth = 7.877604166666667    # 8.7109375
$ws.Columns("O:O").ColumnWidth = 20.022135416666668   # 20.85546875
$ws.Columns("P:P").ColumnWidth = 9.877604166666666    # 10.7109375
$ws.Columns("Q:Q").ColumnWidth = 7.877604166666667    # 8.7109375
$ws.Columns("R:S").ColumnWidth = 20.022135416666668   # 20.85546875

# 4. Row1: wrap text only, on the 3 "Complaints" group headers above the wide columns
$ws.Range("I1").WrapText = $true
$ws.Range("L1").WrapText = $true
$ws.Range("O1").WrapText = $true

# 5. Row3: the "Complaints" header cells also get wrap text (keeps existing bold/border/center)
$ws.Range("I3").WrapText = $true
$ws.Range("L3").WrapText = $true
$ws.Range("O3").WrapText = $true

# 6. Row4: center vertically across the whole row, and wrap text on the wide columns
$ws.Range("A4:S4").VerticalAlignment = -4108
$ws.Range("I4").WrapText = $true
$ws.Range("L4").WrapText = $true
$ws.Range("O4").WrapText = $true
$ws.Range("R4").WrapText = $true
$ws.Range("S4").WrapText = $true

Write-Host "step-done"
